# Auto update stock data
# Updates the "Date_1" (A) and "EBITDA" (B) columns for the first data row
# of each company block (one refresh date per ticker) from 2025/12/15 to
# 2025/12/16, along with refreshed EBITDA figures for that date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> updated EBITDA (column B) value; $null means column B is unchanged
$updates = [ordered]@{
    2  = "5.58"
    8  = "8.06"
    14 = "2.91"
    20 = "12.88"
    26 = "10.83"
    32 = "26.54"
    38 = $null
    44 = "11.07"
    50 = "12.18"
    56 = "31.73"
    62 = "11.98"
    68 = "14.03"
    74 = "16.39"
}

foreach ($row in $updates.Keys) {
    # Column A: refresh date. Force literal text storage (rather than
    # letting Excel auto-convert the string into a date serial number)
    # while keeping the cell's original ("Normal") style.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = "2025/12/16"
    $cellA.Style = "Normal"

    $newB = $updates[$row]
    if ($null -ne $newB) {
        # Column B: EBITDA. Same text-forcing trick so "8.06" etc. stay
        # literal text instead of becoming a numeric value.
        $cellB = $ws.Cells.Item($row, 2)
        $cellB.NumberFormat = "@"
        $cellB.Value = $newB
        $cellB.Style = "Normal"
    }
}
